$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TC002_Output")

$ws.Range("A7").Value = "Dr. Chandrashekara Aithal"
$ws.Range("A8").Value = "Dr. S C Rajendran"
$ws.Range("A10").Value = "Dr. Revanth BN"
$ws.Range("A11").Value = "Dr. Shaik Mohammed Obeidullah"
